$wb = $excel.ActiveWorkbook

# Updated "想去人数" (F column) values for rows in both the "展览" sheet
# and the "全部类型" sheet (which mirrors the same data).
$updates = @{
    3  = 1373
    4  = 157
    7  = 101
    9  = 183
    10 = 135
    11 = 4614
    12 = 6875
    18 = 4139
    19 = 593
    21 = 63
    22 = 2723
    26 = 363
    27 = 371
    29 = 227
    30 = 43
    31 = 1631
    32 = 1025
    33 = 67
    34 = 229
    36 = 548
    39 = 91
    40 = 146
    41 = 647
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
